$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("By Section")

$ws.Range("J3:J8").Formula = "=SUM(B3:I3)"
$ws.Range("M2").Formula = '=A2&" & "&B2&" & "&C2&" & "&D2&" & "&E2&" & "&F2&" & "&G2&" & "&H2&" & "&I2&" & "&J2&" \\"'
$ws.Range("M3").Formula = '=A3&" & "&B3&" & "&C3&" & "&D3&" & "&E3&" & "&F3&" & "&G3&" & "&H3&" & "&I3&" & "&J3&" \\"'
$ws.Range("M4:M8").Formula = '=A4&" & "&B4&" & "&C4&" & "&D4&" & "&E4&" & "&F4&" & "&G4&" & "&H4&" & "&I4&" & "&J4&" \\"'
